$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Templates")

# Row 24 updates
$ws.Range("B24").Value = 201
$ws.Range("C24").Value = 486
$ws.Range("D24").Value = 302
$ws.Range("E24").Value = 537

# Row 57 updates
$ws.Range("B57").Value = 1546
$ws.Range("C57").Value = 280
$ws.Range("D57").Value = 1647
$ws.Range("E57").Value = 322
